$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '91.833.32'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '3.084.23'
$ws.Range('E3').Value = '  -2.09%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '231.56'
$ws.Range('E5').Value = '  -4.12%  '
$ws.Range('D6').Value = '608.55'
$ws.Range('E6').Value = '  -1.63%  '
$ws.Range('D7').Value = '1.07'
$ws.Range('E7').Value = '  -5.48%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.380'
$ws.Range('E8').Value = '  +1.03%  '
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('D10').Value = '3.078.87'
$ws.Range('E10').Value = '  -2.22%  '
$ws.Range('D11').Value = '0.757'
$ws.Range('E11').Value = '  +1.53%  '
$ws.Range('D12').Value = '0.196'
$ws.Range('E12').Value = '  -4.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000240'
$ws.Range('E13').Value = '  -3.46%  '
$ws.Range('D14').Value = '91.730.83'
$ws.Range('E14').Value = '  +0.74%  '
$ws.Range('D15').Value = '33.26'
$ws.Range('E15').Value = '  -5.49%  '
$ws.Range('D16').Value = '5.33'
$ws.Range('E16').Value = '  -5.16%  '
$ws.Range('D17').Value = '3.665.65'
$ws.Range('E17').Value = '  -2.01%  '
$ws.Range('D18').Value = '3.058.35'
$ws.Range('E18').Value = '  -3.12%  '
$ws.Range('D19').Value = '3.74'
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').Value = '14.27'
$ws.Range('E20').Value = '  -5.39%  '
$ws.Range('D21').Value = '5.69'
$ws.Range('E21').Value = '  -5.17%  '
$ws.Range('D22').Value = '431.85'
$ws.Range('E22').Value = '  -5.77%  '
$ws.Range('D23').Value = '8.96'
$ws.Range('E23').Value = '  -2.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000190'
$ws.Range('E24').Value = '  -7.48%  '
$ws.Range('D25').Value = '5.52'
$ws.Range('E25').Value = '  -6.99%  '
$ws.Range('D26').Value = '84.78'
$ws.Range('E26').Value = '  -4.79%  '
$ws.Range('D27').Value = '11.21'
$ws.Range('E27').Value = '  -5.63%  '
$ws.Range('D28').Value = '3.251.75'
$ws.Range('E28').Value = '  -2.06%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = '0.125'
$ws.Range('E30').Value = '  -16.60%  '
$ws.Range('B31').Value = 'Cronos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D31').Value = '0.174'
$ws.Range('E31').Value = '  +3.42%  '
$ws.Range('D32').Value = '0.229'
$ws.Range('E32').Value = '  -4.20%  '
$ws.Range('B33').Value = 'Binance-PegBSC-USD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D33').Value = '1.04'
$ws.Range('E33').Value = '  +48.12%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = '9.01'
$ws.Range('E34').Value = '  -4.79%  '
$ws.Range('D35').Value = '7.71'
$ws.Range('E35').Value = '  +2.50%  '
$ws.Range('D36').Value = '0.155'
$ws.Range('E36').Value = '  -12.00%  '
$ws.Range('D37').Value = '25.16'
$ws.Range('E37').Value = '  -5.50%  '
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').Value = '1.87'
$ws.Range('E39').Value = '  -3.63%  '
$ws.Range('D40').Value = '23.82'
$ws.Range('E40').Value = '  +7.68%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').Value = '463.53'
$ws.Range('E41').Value = '  -6.54%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').Value = '1.26'
$ws.Range('E42').Value = '  -5.43%  '
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').Value = '0.426'
$ws.Range('E43').Value = '  -5.21%  '
$ws.Range('D44').Value = '3.22'
$ws.Range('E44').Value = '  -6.35%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').Value = '160.48'
$ws.Range('E46').Value = '  +2.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.670'
$ws.Range('E47').Value = '  -6.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.80'
$ws.Range('E48').Value = '  -6.68%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').Value = '43.75'
$ws.Range('E49').Value = '  -0.68%  '
$ws.Range('B50').Value = 'ImmutableX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D50').Value = '1.31'
$ws.Range('E50').Value = '  -4.08%  '
$ws.Range('D51').Value = '0.998'
$ws.Range('E51').Value = '  -0.05%  '
